# Insert a new test row ("Фильтрация Courses" / "Can_Filter_Courses") above
# the last two existing test rows, shifting them down by one, matching the
# target diff (old row 23 -> 24, old row 24 -> 25, new row 16 inserted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 16; everything from row 16 down (including the
# two rows at 23/24) shifts down by one row.
$ws.Rows.Item(16).Insert() | Out-Null

# Populate the newly inserted row 16 with the new test data.
$ws.Range("A16").Value = "Фильтрация Courses"
$ws.Range("B16").Value = "Can_Filter_Courses"
$ws.Range("C16").Value = 0

# Match the recorded selection from the diff.
$ws.Range("C16").Select() | Out-Null
